$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.631.56'
$ws.Range('E2').Value = '  +1.55%  '
$ws.Range('D3').Value = '1.602.12'
$ws.Range('E3').Value = '  +1.39%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.44'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('E6').Value = '  +0.92%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '27.83'
$ws.Range('E8').Value = '  +6.42%  '
$ws.Range('E9').Value = '  +1.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0602'
$ws.Range('E10').Value = '  +1.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0910'
$ws.Range('E11').Value = '  +0.48%  '
$ws.Range('D12').Value = '1.832.36'
$ws.Range('E12').Value = '  +1.50%  '
$ws.Range('D13').Value = '1.611.39'
$ws.Range('E13').Value = '  +1.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.544'
$ws.Range('E14').Value = '  +4.10%  '
$ws.Range('D15').Value = '29.640.15'
$ws.Range('E15').Value = '  +1.48%  '
$ws.Range('E16').Value = '  +1.16%  '
$ws.Range('E17').Value = '  +2.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '241.46'
$ws.Range('E18').Value = '  +1.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.70'
$ws.Range('E19').Value = '  +3.45%  '
$ws.Range('E20').Value = '  +0.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('E22').Value = '  +0.38%  '
$ws.Range('E23').Value = '  +1.54%  '
$ws.Range('E24').Value = '  +1.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.79'
$ws.Range('E25').Value = '  +0.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.43'
$ws.Range('E26').Value = '  +1.95%  '
$ws.Range('E27').Value = '  +0.61%  '
$ws.Range('E28').Value = '  +0.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  +2.37%  '
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('E32').Value = '  +0.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.18'
$ws.Range('E33').Value = '  +3.76%  '
$ws.Range('D34').Value = '1.423.13'
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('E35').Value = '  +3.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.89'
$ws.Range('E36').Value = '  +4.31%  '
$ws.Range('E37').Value = '  -1.99%  '
$ws.Range('E38').Value = '  -0.33%  '
$ws.Range('E39').Value = '  +2.87%  '
$ws.Range('E40').Value = '  +3.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '55.47'
$ws.Range('E41').Value = '  +4.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.98'
$ws.Range('E42').Value = '  +0.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0493'
$ws.Range('E43').Value = '  +4.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.815'
$ws.Range('E44').Value = '  +3.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.987'
$ws.Range('E46').Value = '  +17.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '66.08'
$ws.Range('E47').Value = '  +2.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.35'
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D49').Value = '1.742.66'
$ws.Range('E49').Value = '  +1.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '86.66'
$ws.Range('E50').Value = '  +1.65%  '
$ws.Range('D51').Value = '0.0₆0103'
$ws.Range('E51').Value = '  -0.48%  '
